$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Grow the data table: it used to be 3 workers x 3 periods (rows
#    16-24); now it is 3 workers x 4 periods (rows 16-27). Insert 3
#    new rows above the old last row (24) so the trailing "signature"
#    rows (old 29/30) shift down to 32/33, matching the diff.
# ------------------------------------------------------------------
$ws.Range("B24:J26").EntireRow.Insert()

# Copy the formatting of the last "normal" data row (23) onto the
# freshly inserted rows so they pick up styles 15/16/17/16/18/18/19/19/20
# (same look as every other interior row of the table).
$ws.Range("B23:J23").Copy()
$ws.Range("B24:J26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Rewrite the worker/period table (rows 16-27). The sheet keeps
#    the "Tipo Doc" (CC) / Valor Mora (56940) / Salario Basico
#    (1423500) values as before; what changes is which worker each
#    row belongs to and which period it covers - each of the three
#    workers now has 4 periods (2507, 2506, 2505, 2504) instead of 3.
# ------------------------------------------------------------------
$tipoDoc = "CC"
$valorMora = 56940
$salarioBasico = 1423500

$workers = @(
    @{ Doc = "14974528"; Nombre = "BERNABE VERA DIAZ" },
    @{ Doc = "9093657";  Nombre = "ARMANDO ENRIQUE BARRETO PITALUA" },
    @{ Doc = "73594786"; Nombre = "RAFAEL ENRIQUE CONEO GONZALEZ" }
)
$periodos = @("2507", "2506", "2505", "2504")

$row = 16
foreach ($worker in $workers) {
    foreach ($periodo in $periodos) {
        $ws.Cells.Item($row, 2).Value = $tipoDoc
        $ws.Cells.Item($row, 3).Value = $worker.Doc
        $ws.Cells.Item($row, 4).Value = $worker.Nombre
        $ws.Cells.Item($row, 5).Value = $periodo
        $ws.Cells.Item($row, 6).Value = $valorMora
        $ws.Cells.Item($row, 7).Value = $salarioBasico
        $row = $row + 1
    }
}

# ------------------------------------------------------------------
# 3. Update the summary figures: total "Valor Mora" and "Cant.
#    Periodos" (trabajadores count stays at 3, periods count goes
#    from 3 to 4).
# ------------------------------------------------------------------
$ws.Range("E11").Value = 683280
$ws.Range("F13").Value = 4
